$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.292.86"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "1.584.84"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.61"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.807.53"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.580.94"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.42"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "26.291.13"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "207.16"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  -4.06%  "
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.36"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.02"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  +13.40%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "1.283.91"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.615"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.46"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.768"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.39"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "1.719.82"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.84"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.43"
$ws.Range("E51").Value = "  +0.18%  "
